$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header/data rows down by one row to make room for the
# "DELETE THIS ROW WHEN YOU SAVE" note at the top, and make room for the
# extra example rows.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Header row (row 2 after the insert above)
$ws.Range("A2").Value = "Source"
$ws.Range("B2").Value = "Destination"
$ws.Range("C2").Value = "Description"
$ws.Range("D2").Value = "FileExtension"

# Example/help text, entered in the order the author originally wrote them
$ws.Range("D6").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Here are some examples all you have to do is to remove the '' around the text"
$ws.Range("D7").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: '*.txt ' will get all files that end in '.txt'"
$ws.Range("D8").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: 'name* ' will get all files that start with 'name'"
$ws.Range("D9").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: '*name* ' will get all files that contains 'name'"

# Instructional note placed at the very top of the sheet
$ws.Range("A1").Value = "DELETE THIS ROW WHEN YOU SAVE:  The default process is backup, so when you run a backup, it will read from the ""Source"" column and put it into the ""Destination"" column."

$ws.Range("D10").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: 'name*.txt' will get all files that start with 'name' and ends in '.txt'"
$ws.Range("D11").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  As shown by the last example you can combine them to make complex patterns such as the following"
$ws.Range("D12").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Example: 'start * middle * end.txt' will get all files that start with 'start ', and somewhere in the middle has ' middle ', and ends with  ' end.txt'"

# Rows inserted afterwards above the last examples
$ws.Range("D4").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Using a pattern, with the * as any number of characters, you can fetch the desired files seperated by a '/' for every entry"
$ws.Range("D5").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  or you can leave it blank to get ALL the contents of the folder.  This is case insensitive so you don't need to worry about capital letters"
$ws.Range("D3").Value = "REMOVE\REPLACE ME WHEN YOU SAVE:  Need to change the examples because I had to return to the -match comparison operator and not the -like operator"

$ws.PageSetup.Orientation = 1
